$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 540
$ws1.Range("F11").Value = 1603
$ws1.Range("F14").Value = 397
$ws1.Range("F15").Value = 260
$ws1.Range("F21").Value = 182

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 540
$ws4.Range("F12").Value = 1603
$ws4.Range("F15").Value = 397
$ws4.Range("F16").Value = 260
$ws4.Range("F22").Value = 182
